$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing data (a..i) shifts down to rows 2..10
$ws.Rows.Item(1).Insert()

# Append new book title letters j..p in rows 11..17 (before the header text, so
# the shared-string table keeps a..p contiguous and "Book Titles" lands last)
$ws.Range("A11").Value = "j"
$ws.Range("A12").Value = "k"
$ws.Range("A13").Value = "l"
$ws.Range("A14").Value = "m"
$ws.Range("A15").Value = "n"
$ws.Range("A16").Value = "o"
$ws.Range("A17").Value = "p"

# Set the header text
$ws.Range("A1").Value = "Book Titles"

# Style the header cell: Source Sans Pro, size 10, font color FFFFDEDE
$ws.Range("A1").Font.Name = "Source Sans Pro"
$ws.Range("A1").Font.Size = 10
$ws.Range("A1").Font.Color = 14606079
$ws.Range("A1").VerticalAlignment = -4108

# Reset selection to A1 so no stale selection is persisted
$ws.Range("A1").Select()
